$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row2:
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.009821333333333333
$ws.Range("H2").Value = 0.029464
$ws.Range("I2").Value = 0.06297798848338983
$ws.Range("J2").Value = 0.06297798848338984
$ws.Range("M2").Value = 5.079146
$ws.Range("N2").Value = 15.237438
$ws.Range("O2").Value = 0.707738273396773
$ws.Range("P2").Value = 0.7077382733967729
$ws.Range("Q2").Value = 0.04988398591466666
$ws.Range("R2").Value = 0.448955873232
$ws.Range("S2").Value = 0.04457193283123617
$ws.Range("T2").Value = 0.04457193283123617

# Row3:
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.009821333333333333
$ws.Range("H3").Value = 0.029464
$ws.Range("I3").Value = 0.06297798848338983
$ws.Range("J3").Value = 0.06297798848338984
$ws.Range("O3").Value = 0.2087923026002514
$ws.Range("P3").Value = 0.2087923026002514
$ws.Range("Q3").Value = 0.01471644628177778
$ws.Range("R3").Value = 0.132448016536
$ws.Range("S3").Value = 0.01314931922857908
$ws.Range("T3").Value = 0.01314931922857908

# Row4:
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.009821333333333333
$ws.Range("H4").Value = 0.029464
$ws.Range("I4").Value = 0.06297798848338983
$ws.Range("J4").Value = 0.06297798848338984
$ws.Range("M4").Value = 0.5990256666666667
$ws.Range("N4").Value = 1.797077
$ws.Range("O4").Value = 0.08346942400297561
$ws.Range("P4").Value = 0.08346942400297561
$ws.Range("Q4").Value = 0.005883230747555555
$ws.Range("R4").Value = 0.05294907672800001
$ws.Range("S4").Value = 0.00525673642357458
$ws.Range("T4").Value = 0.005256736423574581

# Row5:
$ws.Range("I5").Value = 0.3247949111459754
$ws.Range("J5").Value = 0.3247949111459754
$ws.Range("M5").Value = 5.079146
$ws.Range("N5").Value = 15.237438
$ws.Range("O5").Value = 0.707738273396773
$ws.Range("P5").Value = 0.7077382733967729
$ws.Range("Q5").Value = 0.2572655170946667
$ws.Range("R5").Value = 2.315389653852
$ws.Range("S5").Value = 0.2298697896225109
$ws.Range("T5").Value = 0.2298697896225109

# Row6:
$ws.Range("I6").Value = 0.3247949111459754
$ws.Range("J6").Value = 0.3247949111459754
$ws.Range("O6").Value = 0.2087923026002514
$ws.Range("P6").Value = 0.2087923026002514
$ws.Range("S6").Value = 0.06781467737101227
$ws.Range("T6").Value = 0.06781467737101227

# Row7:
$ws.Range("I7").Value = 0.3247949111459754
$ws.Range("J7").Value = 0.3247949111459754
$ws.Range("M7").Value = 0.5990256666666667
$ws.Range("N7").Value = 1.797077
$ws.Range("O7").Value = 0.08346942400297561
$ws.Range("P7").Value = 0.08346942400297561
$ws.Range("Q7").Value = 0.03034144871755556
$ws.Range("R7").Value = 0.273073038458
$ws.Range("S7").Value = 0.0271104441524522
$ws.Range("T7").Value = 0.02711044415245221

# Row8:
$ws.Range("G8").Value = 0.09547600000000001
$ws.Range("H8").Value = 0.286428
$ws.Range("I8").Value = 0.6122271003706348
$ws.Range("J8").Value = 0.6122271003706349
$ws.Range("M8").Value = 5.079146
$ws.Range("N8").Value = 15.237438
$ws.Range("O8").Value = 0.707738273396773
$ws.Range("P8").Value = 0.7077382733967729
$ws.Range("Q8").Value = 0.484936543496
$ws.Range("R8").Value = 4.364428891464
$ws.Range("S8").Value = 0.4332965509430259
$ws.Range("T8").Value = 0.4332965509430259

# Row9:
$ws.Range("G9").Value = 0.09547600000000001
$ws.Range("H9").Value = 0.286428
$ws.Range("I9").Value = 0.6122271003706348
$ws.Range("J9").Value = 0.6122271003706349
$ws.Range("O9").Value = 0.2087923026002514
$ws.Range("P9").Value = 0.2087923026002514
$ws.Range("Q9").Value = 0.1430627978413333
$ws.Range("R9").Value = 1.287565180572
$ws.Range("S9").Value = 0.1278283060006601
$ws.Range("T9").Value = 0.1278283060006601

# Row10:
$ws.Range("G10").Value = 0.09547600000000001
$ws.Range("H10").Value = 0.286428
$ws.Range("I10").Value = 0.6122271003706348
$ws.Range("J10").Value = 0.6122271003706349
$ws.Range("M10").Value = 0.5990256666666667
$ws.Range("N10").Value = 1.797077
$ws.Range("O10").Value = 0.08346942400297561
$ws.Range("P10").Value = 0.08346942400297561
$ws.Range("Q10").Value = 0.05719257455066667
$ws.Range("R10").Value = 0.514733170956
$ws.Range("S10").Value = 0.05110224342694882
$ws.Range("T10").Value = 0.05110224342694883

Write-Host "Update complete"
